$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1168.0834
$ws.Range("I28").Value = 559.4211
$ws.Range("J28").Value = 3481
$ws.Range("K28").Value = 559.4211
$ws.Range("L28").Value = 3481
$ws.Range("M28").Value = -74.42110000000002
$ws.Range("N28").Value = -4451
# Row 94
$ws.Range("H94").Value = 1280.4
$ws.Range("I94").Value = 1280.4
$ws.Range("K94").Value = 1280.4
$ws.Range("M94").Value = -829.4000000000001
# Row 99
$ws.Range("H99").Value = 2699.6667
$ws.Range("J99").Value = 5088.3335
$ws.Range("L99").Value = 15265.0005
$ws.Range("N99").Value = -18261.0005
# Row 103
$ws.Range("H103").Value = 1424.5883
$ws.Range("J103").Value = 1500
$ws.Range("L103").Value = 4500
$ws.Range("N103").Value = -5672
# Row 106
$ws.Range("H106").Value = 8797.833000000001
$ws.Range("I106").Value = 4966.3335
$ws.Range("K106").Value = 4966.3335
$ws.Range("M106").Value = -4335.3335

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 766.425
$ws.Range("I2").Value = 773.2564
$ws.Range("K2").Value = 773.2564
$ws.Range("M2").Value = -660.2564
# Row 39
$ws.Range("H39").Value = 14428.571
$ws.Range("I39").Value = 14428.571
$ws.Range("K39").Value = 14428.571
$ws.Range("M39").Value = -13908.571
# Row 63
$ws.Range("H63").Value = 3428.25
$ws.Range("I63").Value = 3060.8572
$ws.Range("J63").Value = 6000
$ws.Range("K63").Value = 3060.8572
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = -2374.8572
$ws.Range("N63").Value = -7372
# Row 66
$ws.Range("H66").Value = 3428.25
$ws.Range("I66").Value = 3060.8572
$ws.Range("J66").Value = 6000
$ws.Range("K66").Value = 15304.286
$ws.Range("L66").Value = 30000
$ws.Range("M66").Value = -11872.286
$ws.Range("N66").Value = -36864
# Row 110
$ws.Range("H110").Value = 2367.2104
$ws.Range("I110").Value = 2909.5557
$ws.Range("K110").Value = 2909.5557
$ws.Range("M110").Value = -864.5556999999999
# Row 116
$ws.Range("H116").Value = 766.425
$ws.Range("I116").Value = 773.2564
$ws.Range("K116").Value = 773.2564
$ws.Range("M116").Value = 1520.7436
# Row 130
$ws.Range("H130").Value = 27706.857
$ws.Range("J130").Value = 27706.857
$ws.Range("L130").Value = 27706.857
$ws.Range("N130").Value = -37746.857
# Row 135
$ws.Range("H135").Value = 56570.2
$ws.Range("J135").Value = 56570.2
$ws.Range("L135").Value = 56570.2
$ws.Range("N135").Value = -66710.2

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 766.425
$ws.Range("I3").Value = 773.2564
$ws.Range("K3").Value = 773.2564
$ws.Range("M3").Value = -659.2564
# Row 20
$ws.Range("H20").Value = 3314.125
$ws.Range("I20").Value = 1202.6
$ws.Range("J20").Value = 6833.3335
$ws.Range("K20").Value = 1202.6
$ws.Range("L20").Value = 6833.3335
$ws.Range("M20").Value = -955.5999999999999
$ws.Range("N20").Value = -7327.3335
# Row 86
$ws.Range("H86").Value = 3826.3333
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
# Row 89
$ws.Range("H89").Value = 3826.3333
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
# Row 105
$ws.Range("H105").Value = 3988.682
$ws.Range("I105").Value = 3987.2104
$ws.Range("J105").Value = 3998
$ws.Range("K105").Value = 3987.2104
$ws.Range("L105").Value = 3998
$ws.Range("M105").Value = -2240.2104
$ws.Range("N105").Value = -7492
# Row 107
$ws.Range("H107").Value = 4300
$ws.Range("I107").Value = 3600
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 3600
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = -1680
$ws.Range("N107").Value = -8840

$ws = $wb.Worksheets.Item("CRP")
# Row 68
$ws.Range("H68").Value = 45000
$ws.Range("I68").Value = 30000
$ws.Range("K68").Value = 30000
$ws.Range("M68").Value = -29251
# Row 71
$ws.Range("H71").Value = 45000
$ws.Range("I71").Value = 30000
$ws.Range("K71").Value = 90000
$ws.Range("M71").Value = -86256
# Row 86
$ws.Range("H86").Value = 9869.125
$ws.Range("I86").Value = 5810.8
$ws.Range("J86").Value = 16633
$ws.Range("K86").Value = 5810.8
$ws.Range("L86").Value = 16633
$ws.Range("M86").Value = -4687.8
$ws.Range("N86").Value = -18879
# Row 89
$ws.Range("H89").Value = 9869.125
$ws.Range("I89").Value = 5810.8
$ws.Range("J89").Value = 16633
$ws.Range("K89").Value = 29054
$ws.Range("L89").Value = 83165
$ws.Range("M89").Value = -23438
$ws.Range("N89").Value = -94397
# Row 105
$ws.Range("H105").Value = 1250
$ws.Range("I105").Value = 1000
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 1000
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 747
$ws.Range("N105").Value = -4994

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 1850.3334
$ws.Range("I14").Value = 1850.3334
$ws.Range("K14").Value = 5551.0002
$ws.Range("M14").Value = -5378.0002
# Row 107
$ws.Range("H107").Value = 455632.5
$ws.Range("J107").Value = 556767.5
$ws.Range("L107").Value = 1670302.5
$ws.Range("N107").Value = -1674142.5
# Row 111
$ws.Range("H111").Value = 4714.5557
$ws.Range("I111").Value = 1739.3334
$ws.Range("J111").Value = 10665
$ws.Range("K111").Value = 5218.0002
$ws.Range("L111").Value = 31995
$ws.Range("M111").Value = -2151.0002
$ws.Range("N111").Value = -38129
# Row 131
$ws.Range("H131").Value = 3232823.5
$ws.Range("I131").Value = 995.2
$ws.Range("K131").Value = 2985.6
$ws.Range("M131").Value = 2054.4
# Row 136
$ws.Range("H136").Value = 3369.75
$ws.Range("I136").Value = 3369.75
$ws.Range("K136").Value = 10109.25
$ws.Range("M136").Value = -5009.25

$ws = $wb.Worksheets.Item("GSM")
# Row 114
$ws.Range("H114").Value = 80000
$ws.Range("J114").Value = 80000
$ws.Range("L114").Value = 80000
$ws.Range("N114").Value = -88678

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 7334.0884
$ws.Range("I61").Value = 6859.769
$ws.Range("J61").Value = 8875.625
$ws.Range("K61").Value = 6859.769
$ws.Range("L61").Value = 8875.625
$ws.Range("M61").Value = -6657.769
$ws.Range("N61").Value = -9279.625
# Row 113
$ws.Range("H113").Value = 7334.0884
$ws.Range("I113").Value = 6859.769
$ws.Range("J113").Value = 8875.625
$ws.Range("K113").Value = 6859.769
$ws.Range("L113").Value = 8875.625
$ws.Range("M113").Value = -4689.769
$ws.Range("N113").Value = -13215.625
# Row 132
$ws.Range("H132").Value = 82152.07000000001
$ws.Range("I132").Value = 87631.21000000001
$ws.Range("K132").Value = 262893.63
$ws.Range("M132").Value = -260363.63
# Row 133
$ws.Range("H133").Value = 99999
$ws.Range("J133").Value = 99999
$ws.Range("L133").Value = 99999
$ws.Range("N133").Value = -105059
# Row 136
$ws.Range("H136").Value = 5178.875
$ws.Range("I136").Value = 3549.75
$ws.Range("K136").Value = 10649.25
$ws.Range("M136").Value = -8099.25

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 97486
$ws.Range("I62").Value = 5860.75
$ws.Range("J62").Value = 149843.28
$ws.Range("K62").Value = 5860.75
$ws.Range("L62").Value = 149843.28
$ws.Range("M62").Value = -5236.75
$ws.Range("N62").Value = -151091.28
# Row 65
$ws.Range("H65").Value = 97486
$ws.Range("I65").Value = 5860.75
$ws.Range("J65").Value = 149843.28
$ws.Range("K65").Value = 29303.75
$ws.Range("L65").Value = 749216.4
$ws.Range("M65").Value = -26183.75
$ws.Range("N65").Value = -755456.4
# Row 81
$ws.Range("H81").Value = 1350
$ws.Range("I81").Value = 1350
$ws.Range("K81").Value = 2700
$ws.Range("M81").Value = -1639
# Row 84
$ws.Range("H84").Value = 1350
$ws.Range("I84").Value = 1350
$ws.Range("K84").Value = 13500
$ws.Range("M84").Value = -8196
# Row 117
$ws.Range("H117").Value = 20682
$ws.Range("J117").Value = 20682
$ws.Range("L117").Value = 20682
$ws.Range("N117").Value = -29860
